$wb = $excel.ActiveWorkbook

$wsReq = $wb.Worksheets.Item("Requerimientos")
$wsMustHave = $wb.Worksheets.Item("Must Have")

# Update B2, B10, B11 on Requerimientos sheet from "p" to "x"
$wsReq.Range("B2").Value = "x"
$wsReq.Range("B10").Value = "x"
$wsReq.Range("B11").Value = "x"

# Set selections / active cells per the diff.
# "Must Have" was the active sheet before; move its selection to F2.
$wsMustHave.Activate()
$wsMustHave.Range("F2").Select()

# "Requerimientos" becomes the active sheet, with B16 selected.
$wsReq.Activate()
$wsReq.Range("B16").Select()

$wb.Save()
